$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Jurassic" failed-test count from 202 to 84
$ws.Range("B5").Value = 84

# Move the active selection to B6, mirroring the author's final cursor position
$ws.Range("B6").Select()
